# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The whole column was bumped by one day (serial 45171 -> 45172),
# reflecting an automatic re-export of the workbook on the next day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45172
